$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.077.49'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '1.891.39'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.67'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5173'
$ws.Range('E7').Value = '  +2.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3754'
$ws.Range('E8').Value = '  +2.45%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07208'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.10'
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8965'
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07654'
$ws.Range('E12').Value = '  +1.70%  '
$ws.Range('D13').Value = '1.894.08'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.28'
$ws.Range('E14').Value = '  -0.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.229'
$ws.Range('E15').Value = '  -0.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9998'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008513'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.38'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9996'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '27.135.21'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('D22').Value = '2.127.77'
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.59'
$ws.Range('E23').Value = '  +1.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.404'
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.280'
$ws.Range('E25').Value = '  +9.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.97'
$ws.Range('E26').Value = '  -1.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.732'
$ws.Range('E27').Value = '  -3.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.04'
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.28'
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.966'
$ws.Range('E30').Value = '  +5.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.771'
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09189'
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05044'
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.236'
$ws.Range('E34').Value = '  +6.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7746'
$ws.Range('E35').Value = '  +2.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.981'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.277'
$ws.Range('E37').Value = '  +1.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.586'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5598'
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01987'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.074'
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.986'
$ws.Range('E42').Value = '  +5.03%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '119.24'
$ws.Range('E43').Value = '  +2.82%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.626'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('E45').Value = '  +2.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4817'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.17'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9994'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E49').Value = '  +1.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.37'
$ws.Range('E50').Value = '  +1.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.80'
$ws.Range('E51').Value = '  +0.85%  '
